$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07024999999999999
$ws.Range("H2").Value = 0.21075
$ws.Range("I2").Value = 0.005236595731231519
$ws.Range("J2").Value = 0.005236595731231519
$ws.Range("M2").Value = 24.75542533333333
$ws.Range("N2").Value = 74.26627599999999
$ws.Range("O2").Value = 0.7762421087066456
$ws.Range("P2").Value = 0.7762421087066456
$ws.Range("Q2").Value = 1.739068629666666
$ws.Range("R2").Value = 15.651617667
$ws.Range("S2").Value = 0.004064866112855373
$ws.Range("T2").Value = 0.004064866112855373

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07024999999999999
$ws.Range("H3").Value = 0.21075
$ws.Range("I3").Value = 0.005236595731231519
$ws.Range("J3").Value = 0.005236595731231519
$ws.Range("M3").Value = 3.818542
$ws.Range("N3").Value = 11.455626
$ws.Range("O3").Value = 0.1197358984688377
$ws.Range("P3").Value = 0.1197358984688377
$ws.Range("Q3").Value = 0.2682525755
$ws.Range("R3").Value = 2.4142731795
$ws.Range("S3").Value = 0.000627008494797086
$ws.Range("T3").Value = 0.000627008494797086

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07024999999999999
$ws.Range("H4").Value = 0.21075
$ws.Range("I4").Value = 0.005236595731231519
$ws.Range("J4").Value = 0.005236595731231519
$ws.Range("M4").Value = 3.317404
$ws.Range("N4").Value = 9.952212
$ws.Range("O4").Value = 0.1040219928245168
$ws.Range("P4").Value = 0.1040219928245168
$ws.Range("Q4").Value = 0.2330476309999999
$ws.Range("R4").Value = 2.097428679
$ws.Range("S4").Value = 0.0005447211235790603
$ws.Range("T4").Value = 0.0005447211235790603

$ws.Range("I5").Value = 0.4287876899474159
$ws.Range("J5").Value = 0.4287876899474159
$ws.Range("M5").Value = 24.75542533333333
$ws.Range("N5").Value = 74.26627599999999
$ws.Range("O5").Value = 0.7762421087066456
$ws.Range("P5").Value = 0.7762421087066456
$ws.Range("Q5").Value = 142.3999977556831
$ws.Range("R5").Value = 1281.599979801148
$ws.Range("S5").Value = 0.3328430606322335
$ws.Range("T5").Value = 0.3328430606322335

$ws.Range("I6").Value = 0.4287876899474159
$ws.Range("J6").Value = 0.4287876899474159
$ws.Range("M6").Value = 3.818542
$ws.Range("N6").Value = 11.455626
$ws.Range("O6").Value = 0.1197358984688377
$ws.Range("P6").Value = 0.1197358984688377
$ws.Range("Q6").Value = 21.96530113735534
$ws.Range("R6").Value = 197.687710236198
$ws.Range("S6").Value = 0.05134127930823124
$ws.Range("T6").Value = 0.05134127930823124

$ws.Range("I7").Value = 0.4287876899474159
$ws.Range("J7").Value = 0.4287876899474159
$ws.Range("M7").Value = 3.317404
$ws.Range("N7").Value = 9.952212
$ws.Range("O7").Value = 0.1040219928245168
$ws.Range("P7").Value = 0.1040219928245168
$ws.Range("Q7").Value = 19.08261788249733
$ws.Range("R7").Value = 171.743560942476
$ws.Range("S7").Value = 0.04460335000695122
$ws.Range("T7").Value = 0.04460335000695122

$ws.Range("G8").Value = 7.592679666666666
$ws.Range("H8").Value = 22.778039
$ws.Range("I8").Value = 0.5659757143213526
$ws.Range("J8").Value = 0.5659757143213525
$ws.Range("M8").Value = 24.75542533333333
$ws.Range("N8").Value = 74.26627599999999
$ws.Range("O8").Value = 0.7762421087066456
$ws.Range("P8").Value = 0.7762421087066456
$ws.Range("Q8").Value = 187.9600145680849
$ws.Range("R8").Value = 1691.640131112764
$ws.Range("S8").Value = 0.4393341819615568
$ws.Range("T8").Value = 0.4393341819615567

$ws.Range("G9").Value = 7.592679666666666
$ws.Range("H9").Value = 22.778039
$ws.Range("I9").Value = 0.5659757143213526
$ws.Range("J9").Value = 0.5659757143213525
$ws.Range("M9").Value = 3.818542
$ws.Range("N9").Value = 11.455626
$ws.Range("O9").Value = 0.1197358984688377
$ws.Range("P9").Value = 0.1197358984688377
$ws.Range("Q9").Value = 28.99296619971267
$ws.Range("R9").Value = 260.936695797414
$ws.Range("S9").Value = 0.06776761066580936
$ws.Range("T9").Value = 0.06776761066580934

$ws.Range("G10").Value = 7.592679666666666
$ws.Range("H10").Value = 22.778039
$ws.Range("I10").Value = 0.5659757143213526
$ws.Range("J10").Value = 0.5659757143213525
$ws.Range("M10").Value = 3.317404
$ws.Range("N10").Value = 9.952212
$ws.Range("O10").Value = 0.1040219928245168
$ws.Range("P10").Value = 0.1040219928245168
$ws.Range("Q10").Value = 25.18798589691866
$ws.Range("R10").Value = 226.691873072268
$ws.Range("S10").Value = 0.0588739216939865
$ws.Range("T10").Value = 0.05887392169398649
